$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LP Github repos")

$data = @(
    @("xAI-Console-UI", "Yes (NPM)", "No", "No", "No", "No", "No", "No"),
    @("xAI-Inference-server", "Yes (NPM)", "No", "No", "No", "No", "No", "No"),
    @("xAI-Labelling-UI", "Yes (NPM and Yarn)", "No", "No", "No", "No", "No", "No"),
    @("xapi-js-client", "Yes (NPM)", "Renovate and Dependabot", "No", "No", "No", "No", "No"),
    @("yum-api-qatests", "Yes (NPM)", "No", "No", "No", "No", "No", "No")
)

$startRow = 17
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le $rowData.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}
